$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D to make room for the new quarterly
# data (Dec-2018 and Sep-2018), shifting the existing quarters right.
$ws.Columns("D:E").Insert()

# Copy number formats/styles from column F (the old column D, now shifted)
# into the two newly inserted columns so the new cells pick up the same
# style indices (date format on the header rows, #,##0 on the data rows).
$ws.Range("F7:F102").Copy()
$ws.Range("D7:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Rows 37 and 79 are section-header rows that never had data in columns
# C:K, so the format-copy above should not leave stray empty cells there.
$ws.Range("D37:E37").Clear()
$ws.Range("D79:E79").Clear()

# Rows 36 and 78 are fully blank spacer rows with no cells at all; undo any
# stray cell the format-copy may have introduced there.
$ws.Range("D36:E36").Clear()
$ws.Range("D78:E78").Clear()

# Give the two new columns a sensible width matching their neighbors.
$ws.Columns("D:E").ColumnWidth = 13.75

$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 338000
$ws.Range("E8").Value = 266700
$ws.Range("D9").Value = 115200
$ws.Range("E9").Value = 98100
$ws.Range("D10").Value = 222800
$ws.Range("E10").Value = 168600
$ws.Range("D12").Value = 275300
$ws.Range("E12").Value = 50100
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("D17").Value = 502600
$ws.Range("E17").Value = 252800
$ws.Range("D18").Value = -164600
$ws.Range("E18").Value = 13900
$ws.Range("D20").Value = -4100
$ws.Range("E20").Value = 35800
$ws.Range("D21").Value = -160100
$ws.Range("E21").Value = 57500
$ws.Range("D22").Value = 8200
$ws.Range("E22").Value = 4900
$ws.Range("D23").Value = -176900
$ws.Range("E23").Value = 44800
$ws.Range("D24").Value = 3200
$ws.Range("E24").Value = -1800
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = -180100
$ws.Range("E26").Value = 46600
$ws.Range("D27").Value = -180100
$ws.Range("E27").Value = 46600
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = 400
$ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 4100
$ws.Range("E32").Value = -35800
$ws.Range("D33").Value = -179700
$ws.Range("E33").Value = 46600
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = -179700
$ws.Range("E35").Value = 46600
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 1137000
$ws.Range("E41").Value = 367400
$ws.Range("D42").Value = 248600
$ws.Range("E42").Value = 301300
$ws.Range("D43").Value = 226700
$ws.Range("E43").Value = 169800
$ws.Range("D44").Value = 70700
$ws.Range("E44").Value = 55200
$ws.Range("D45").Value = 16500
$ws.Range("E45").Value = 20400
$ws.Range("D46").Value = 1699500
$ws.Range("E46").Value = 914100
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("D48").Value = 183100
$ws.Range("E48").Value = 169600
$ws.Range("D49").Value = 18700
$ws.Range("E49").Value = 18400
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 14700
$ws.Range("E52").Value = 12200
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 1916000
$ws.Range("E54").Value = 1114300
$ws.Range("D57").Value = 75500
$ws.Range("E57").Value = 57700
$ws.Range("D58").Value = "NA"
$ws.Range("E58").Value = "NA"
$ws.Range("D59").Value = 146900
$ws.Range("E59").Value = 136200
$ws.Range("D60").Value = 222400
$ws.Range("E60").Value = 193900
$ws.Range("D61").Value = 1017600
$ws.Range("E61").Value = 346100
$ws.Range("D62").Value = 12700
$ws.Range("E62").Value = 12100
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 1252700
$ws.Range("E66").Value = 552100
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = -798900
$ws.Range("E72").Value = -619200
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 663300
$ws.Range("E76").Value = 562200
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = -179700
$ws.Range("E81").Value = 46600
$ws.Range("D83").Value = 8600
$ws.Range("E83").Value = 7800
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 33700
$ws.Range("E89").Value = 56800
$ws.Range("D91").Value = -17700
$ws.Range("E91").Value = -23800
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = 30900
$ws.Range("E94").Value = 4600
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = 701600
$ws.Range("E100").Value = 5400
$ws.Range("D101").Value = 3200
$ws.Range("E101").Value = 400
$ws.Range("D102").Value = 769400
$ws.Range("E102").Value = 67200

Write-Host "Applied DXCM quarterly update"
